# Updates cryptos list values per the latest data refresh (GitHub Actions run).
# Numeric-looking text values (e.g. "221.51") are written with a leading apostrophe so Excel
# keeps them as text (matching the sheet's existing text-formatted Price/Volume columns)
# instead of silently converting them to numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Bitcoin
$ws.Range('D2').Value = '90.552.28'
$ws.Range('E2').Value = '  +3.65%  '

# Row 3: Ethereum
$ws.Range('D3').Value = '3.203.49'
$ws.Range('E3').Value = '  +1.62%  '

# Row 5: Solana
$ws.Range('D5').Value = '''221.51'
$ws.Range('E5').Value = '  +7.53%  '

# Row 6: BNB
$ws.Range('D6').Value = '''638.63'
$ws.Range('E6').Value = '  +5.40%  '

# Row 7: Dogecoin
$ws.Range('D7').Value = '''0.399'
$ws.Range('E7').Value = '  +7.41%  '

# Row 8: XRP
$ws.Range('D8').Value = '''0.704'
$ws.Range('E8').Value = '  +6.74%  '

# Row 9: USDC
$ws.Range('D9').Value = '''1.00'
$ws.Range('E9').Value = '  +0.00%  '

# Row 10: LidoStakedEther
$ws.Range('D10').Value = '3.199.94'

# Row 11: Cardano
$ws.Range('D11').Value = '''0.575'
$ws.Range('E11').Value = '  +8.53%  '

# Row 12: TRON
$ws.Range('E12').Value = '  +3.33%  '

# Row 13: ShibaInu
$ws.Range('D13').Value = '''0.0000258'

# Row 14: Toncoin
$ws.Range('D14').Value = '''5.45'
$ws.Range('E14').Value = '  +4.65%  '

# Row 15: Avalanche
$ws.Range('D15').Value = '''33.43'
$ws.Range('E15').Value = '  +4.90%  '

# Row 16: WrappedBTC
$ws.Range('D16').Value = '90.328.58'
$ws.Range('E16').Value = '  +3.57%  '

# Row 17: WrappedliquidstakedEther2.0
$ws.Range('D17').Value = '3.795.42'
$ws.Range('E17').Value = '  +1.42%  '

# Row 18: WrappedEther
$ws.Range('D18').Value = '3.223.72'
$ws.Range('E18').Value = '  +2.26%  '

# Row 19: SuiNetwork
$ws.Range('B19').Value = 'SuiNetwork'
$ws.Range('C19').Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range('D19').Value = '''3.35'
$ws.Range('E19').Value = '  +10.28%  '

# Row 20: PEPE
$ws.Range('B20').Value = 'PEPE'
$ws.Range('C20').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D20').Value = '''0.0000225'
$ws.Range('E20').Value = '  +76.73%  '

# Row 21: Chainlink
$ws.Range('D21').Value = '''13.49'
$ws.Range('E21').Value = '  +1.41%  '

# Row 22: BitcoinCash
$ws.Range('D22').Value = '''438.58'
$ws.Range('E22').Value = '  +6.61%  '

# Row 23: Uniswap
$ws.Range('D23').Value = '''8.64'
$ws.Range('E23').Value = '  +2.96%  '

# Row 24: Polkadot
$ws.Range('D24').Value = '''5.07'
$ws.Range('E24').Value = '  +1.05%  '

# Row 25: NEARProtocol
$ws.Range('D25').Value = '''5.33'
$ws.Range('E25').Value = '  +4.13%  '

# Row 26: Aptos
$ws.Range('D26').Value = '''11.90'
$ws.Range('E26').Value = '  +1.18%  '

# Row 27: Litecoin
$ws.Range('D27').Value = '''81.30'
$ws.Range('E27').Value = '  +11.49%  '

# Row 28: WrappedeETH
$ws.Range('D28').Value = '3.379.21'
$ws.Range('E28').Value = '  +1.49%  '

# Row 29: Dai
$ws.Range('D29').Value = '''1.00'
$ws.Range('E29').Value = '  +0.04%  '

# Row 30: Cronos
$ws.Range('E30').Value = '  +1.23%  '

# Row 31: Binance-PegBSC-USD
$ws.Range('D31').Value = '''0.998'
$ws.Range('E31').Value = '  -0.11%  '

# Row 32: dogwifhat
$ws.Range('D32').Value = '''4.19'
$ws.Range('E32').Value = '  +42.20%  '

# Row 33: InternetComputer(DFINITY)
$ws.Range('D33').Value = '''8.44'
$ws.Range('E33').Value = '  +3.96%  '

# Row 34: Bittensor
$ws.Range('D34').Value = '''539.69'
$ws.Range('E34').Value = '  +0.07%  '

# Row 35: RenderToken
$ws.Range('D35').Value = '''7.10'
$ws.Range('E35').Value = '  +7.27%  '

# Row 36: PancakeSwap
$ws.Range('D36').Value = '''1.92'
$ws.Range('E36').Value = '  +4.64%  '

# Row 37: Fetch.AI
$ws.Range('D37').Value = '''1.29'
$ws.Range('E37').Value = '  +1.21%  '

# Row 38: EthereumClassic
$ws.Range('D38').Value = '''22.52'
$ws.Range('E38').Value = '  +4.28%  '

# Row 39: WhiteBITCoin
$ws.Range('D39').Value = '''22.38'
$ws.Range('E39').Value = '  +2.52%  '

# Row 40: FirstDigitalUSD
$ws.Range('E40').Value = '  +0.10%  '

# Row 41: Kaspa
$ws.Range('E41').Value = '  -3.05%  '

# Row 42: Stacks
$ws.Range('D42').Value = '''1.94'
$ws.Range('E42').Value = '  +2.95%  '

# Row 44: PolygonEcosystemToken
$ws.Range('D44').Value = '''0.374'
$ws.Range('E44').Value = '  +2.66%  '

# Row 45: Monero
$ws.Range('D45').Value = '''146.13'
$ws.Range('E45').Value = '  -1.59%  '

# Row 46: OKB
$ws.Range('D46').Value = '''44.84'
$ws.Range('E46').Value = '  +4.06%  '

# Row 47: Aave
$ws.Range('D47').Value = '''173.44'
$ws.Range('E47').Value = '  +1.28%  '

# Row 48: Stellar
$ws.Range('B48').Value = 'Stellar'
$ws.Range('C48').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D48').Value = '''0.125'
$ws.Range('E48').Value = '  +1.54%  '

# Row 49: Mantle
$ws.Range('B49').Value = 'Mantle'
$ws.Range('C49').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D49').Value = '''0.752'
$ws.Range('E49').Value = '  +9.14%  '

# Row 50: ImmutableX
$ws.Range('B50').Value = 'ImmutableX'
$ws.Range('C50').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D50').Value = '''1.24'
$ws.Range('E50').Value = '  +2.85%  '

# Row 51: ARBITRUM
$ws.Range('B51').Value = 'ARBITRUM'
$ws.Range('C51').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D51').Value = '''0.623'
$ws.Range('E51').Value = '  +7.89%  '
